$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("Indicar la calificación IRCA de los cuerpos de agua que tienen un nivel de riesgo SIN RIESGO, en caso de no haber ninguno devolver NA.")
Write-Output "found=$found start=$($rng.Start) end=$($rng.End)"

$insPoint = $rng.Start + 92
$target = $d.Range($insPoint, $insPoint)
$target.InsertAfter(" separados por espacio")
Write-Output "step1 done"

$dummy = $d.Paragraphs.Count
Write-Output "dummy=$dummy"

$insPoint2 = $insPoint + 23
$target2 = $d.Range($insPoint2, $insPoint2)
$target2.InsertAfter("ZZZ")
Write-Output "step2 done"

$check = $d.Range(2204, 2470)
Write-Output "check text: [$($check.Text)]"
